# Katalon Studio IVY BIMBO
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header columns P1:T1 (copy style from O1 which already carries the header format)
$ws.Range("P1").Value = "IEPS"
$ws.Range("Q1").Value = "IVA"
$ws.Range("R1").Value = "NOTAX"
$ws.Range("S1").Value = "Item Disc"
$ws.Range("T1").Value = "Category Disc"
$ws.Range("O1").Copy()
$ws.Range("P1:T1").PasteSpecial(-4122)

# Row 76: add TOTAL label
$ws.Range("B76").Value = "TOTAL           `$"

# Row 79: extend TC id list and add new invoice column
$ws.Range("A79").Value = "TC_ID_103,104"
$ws.Range("C79").Value = "ExtraInvoice-1"

# Row 82: extend TC id list
$ws.Range("A82").Value = "TC_32,TC_31"

# Row 84: replace TC id and add new numeric columns
$ws.Range("A84").Value = ",TC_39"
$ws.Range("P84").Value = 8
$ws.Range("R84").Value = 0
$ws.Range("S84").Value = 5
$ws.Range("T84").Value = 10

# Row 85: replace TC id and add new numeric columns
$ws.Range("A85").Value = "TC_43"
$ws.Range("Q85").Value = 16
$ws.Range("S85").Value = 5
$ws.Range("T85").Value = 20

# Row 89: update numeric value
$ws.Range("B89").Value = 1234567890

# Update the view to match the post-edit selection/scroll position
$ws.Range("B89").Select()
$excel.ActiveWindow.ScrollRow = 70
